# "tramas y tipos de dispositivos"
# Adds a new worksheet "Enlace" after "Unidades" with four small
# lookup tables describing switch port commands, switch device types,
# Ethernet frame types and switching methods.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Zoom / view tweaks on the existing "Unidades" sheet -----------------
$ws1.Activate()
$excel.ActiveWindow.Zoom = 130

# --- Create the new sheet, placed right after "Unidades" -----------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Enlace"

# Column widths (approximate the original bestFit widths)
$ws2.Columns.Item(1).ColumnWidth = 11.736979166666666
$ws2.Columns.Item(2).ColumnWidth = 35.022135416666664
$ws2.Columns.Item(3).ColumnWidth = 3.5924479166666665
$ws2.Columns.Item(4).ColumnWidth = 23.877604166666668
$ws2.Columns.Item(5).ColumnWidth = 43.307291666666664

# --- Table 1 (Tabla5): Comandos / Parametros ------------------------------
$ws2.Range("A1").Value = "Comandos"
$ws2.Range("B1").Value = "Parametros"
$ws2.Range("A2").Value = "MDIX"
$ws2.Range("B2").Value = "Deteccion de interfaz cruzada (AUTO)"
$ws2.Range("A3").Value = "DUPLEX"
$ws2.Range("B3").Value = "envio y recepcion (HALF/FULL/AUTO)"
$ws2.Range("A4").Value = "SPEED"
$ws2.Range("B4").Value = "Ancho de banda (10/100/1000/AUTO)"
$ws2.Range("A5").Value = "CDP"
$ws2.Range("B5").Value = "Descubrimiento Vecino Cisco (ENABLE)"
$ws2.Range("A6").Value = "LLDP"
$ws2.Range("B6").Value = "Descubrimiento (TRANSMIT/RECEIVE)"

$rng1 = $ws2.Range("A1:B6")
$tbl1 = $ws2.ListObjects.Add(1, $rng1, 0, 1)
$tbl1.Name = "Tabla5"
$tbl1.TableStyle = "TableStyleMedium16"

# --- Table 2 (Tabla6): Switches / Descripcion -----------------------------
$ws2.Range("A8").Value = "Switches"
$ws2.Range("B8").Value = "Descripcion"
$ws2.Range("A9").Value = "FIJA"
$ws2.Range("B9").Value = "Cantidad de interfaces predeterminada, sin capacidad de expansion"
$ws2.Range("A10").Value = "EXPANDIBLE"
$ws2.Range("B10").Value = "Poseen buses de expansion que permiten agregar nuevas interfaces."
$ws2.Range("A11").Value = "APILABLE"
$ws2.Range("B11").Value = "Permiten la conexión entre varios, para funcionar como uno de mayor capacidad."

$rng2 = $ws2.Range("A8:B11")
$tbl2 = $ws2.ListObjects.Add(1, $rng2, 0, 1)
$tbl2.Name = "Tabla6"
$tbl2.TableStyle = "TableStyleMedium16"

# --- Table 3 (Tabla7): Tipos de Tramas / Detalles -------------------------
$ws2.Range("D1").Value = "Tipos de Tramas"
$ws2.Range("E1").Value = "Detalles"
$ws2.Range("D2").Value = "Runt"
$ws2.Range("E2").Value = "Tramas menores a 64B, productos de una colision"
$ws2.Range("D3").Value = "Giant"
$ws2.Range("E3").Value = "Tramas mayores a 1,5KB, por fallos en la interfaz"
$ws2.Range("D4").Value = "Throttle"
$ws2.Range("E4").Value = "Tramas descartadas, por desbordamiento de buffer"
$ws2.Range("D5").Value = "CRC"
$ws2.Range("E5").Value = "Tramas Corruptas, consecuencia de EMI o RFI"
$ws2.Range("D6").Value = "Jumbo"
$ws2.Range("E6").Value = "Tramas de 9KB, establecidas por configuracion"

$rng3 = $ws2.Range("D1:E6")
$tbl3 = $ws2.ListObjects.Add(1, $rng3, 0, 1)
$tbl3.Name = "Tabla7"
$tbl3.TableStyle = "TableStyleMedium16"

# --- Table 4 (Tabla8): Metodo / Funcionamiento ----------------------------
$ws2.Range("D8").Value = "Metodo"
$ws2.Range("E8").Value = "Funcionamiento"
$ws2.Range("D9").Value = "Cut-Through"
$ws2.Range("E9").Value = "Reenvian la trama apenas la reciben, sin verificar la misma."
$ws2.Range("D10").Value = "Store-N-Forward"
$ws2.Range("E10").Value = "Almacenan la trama en un buffer y la reenvia despues de verificarla"
$ws2.Range("D11").Value = " Fragment-Free"
$ws2.Range("E11").Value = "Leen los primeros 64B y luego realiza el envio de la misma."

$rng4 = $ws2.Range("D8:E11")
$tbl4 = $ws2.ListObjects.Add(1, $rng4, 0, 1)
$tbl4.Name = "Tabla8"
$tbl4.TableStyle = "TableStyleMedium16"

# --- Row heights + wrap/vertical alignment for the taller description rows -
$ws2.Rows.Item(9).RowHeight = 30
$ws2.Rows.Item(10).RowHeight = 30
$ws2.Rows.Item(11).RowHeight = 30

$ws2.Range("A9:A11").VerticalAlignment = -4108

$ws2.Range("B9:B11").VerticalAlignment = -4108
$ws2.Range("B9:B11").WrapText = $true

$ws2.Range("D9:E9").HorizontalAlignment = -4131
$ws2.Range("D9:E9").VerticalAlignment = -4108
$ws2.Range("D9:E9").WrapText = $true

$ws2.Range("D10:E11").VerticalAlignment = -4108
$ws2.Range("D10:E11").WrapText = $true

# --- Selection / activation so the new sheet matches the saved view -------
$ws2.Range("D6").Select()
$ws2.Activate()
$excel.ActiveWindow.Zoom = 145
